$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1692913385826772
$ws.Range("C2").Value = 0.610236220472441
$ws.Range("J2").Value = 0.003937007874015748
$ws.Range("P2").Value = 0.1181102362204724
$ws.Range("S2").Value = 0.09842519685039371
$ws.Range("B3").Value = 0.006211180124223602
$ws.Range("C3").Value = 0.03726708074534162
$ws.Range("J3").Value = 0.03726708074534162
$ws.Range("P3").Value = 0.7018633540372671
$ws.Range("S3").Value = 0.2173913043478261
$ws.Range("P4").Value = 0.7872340425531915
$ws.Range("S4").Value = 0.2127659574468085
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.02985074626865672
$ws.Range("D6").Value = 0.01492537313432836
$ws.Range("E6").Value = 0.007462686567164179
$ws.Range("F6").Value = 0.05223880597014925
$ws.Range("J6").Value = 0.2761194029850746
$ws.Range("O6").Value = 0.003731343283582089
$ws.Range("Q6").Value = 0.1716417910447761
$ws.Range("R6").Value = 0.06716417910447761
$ws.Range("S6").Value = 0.376865671641791
$ws.Range("B7").Value = 0.07199999999999999
$ws.Range("D7").Value = 0.016
$ws.Range("E7").Value = 0.004
$ws.Range("F7").Value = 0.08
$ws.Range("J7").Value = 0.116
$ws.Range("O7").Value = 0.016
$ws.Range("Q7").Value = 0.196
$ws.Range("R7").Value = 0.064
$ws.Range("S7").Value = 0.436
$ws.Range("B8").Value = 0.08749999999999999
$ws.Range("D8").Value = 0.01875
$ws.Range("E8").Value = 0.004166666666666667
$ws.Range("F8").Value = 0.05
$ws.Range("J8").Value = 0.09375
$ws.Range("O8").Value = 0.01666666666666667
$ws.Range("Q8").Value = 0.1916666666666667
$ws.Range("R8").Value = 0.08333333333333333
$ws.Range("S8").Value = 0.4541666666666667
$ws.Range("B9").Value = 0.07547169811320754
$ws.Range("D9").Value = 0.02830188679245283
$ws.Range("F9").Value = 0.08490566037735849
$ws.Range("J9").Value = 0.0880503144654088
$ws.Range("O9").Value = 0.006289308176100629
$ws.Range("Q9").Value = 0.2012578616352201
$ws.Range("R9").Value = 0.07861635220125786
$ws.Range("S9").Value = 0.4371069182389937
$ws.Range("B10").Value = 0.08878504672897196
$ws.Range("D10").Value = 0.01713395638629283
$ws.Range("E10").Value = 0.000778816199376947
$ws.Range("F10").Value = 0.08099688473520249
$ws.Range("J10").Value = 0.1020249221183801
$ws.Range("O10").Value = 0.008566978193146417
$ws.Range("Q10").Value = 0.221183800623053
$ws.Range("R10").Value = 0.07398753894080996
$ws.Range("S10").Value = 0.4065420560747663
$ws.Range("G11").Value = 0.1459627329192547
$ws.Range("J11").Value = 0.07142857142857142
$ws.Range("K11").Value = 0.1832298136645963
$ws.Range("L11").Value = 0.5962732919254659
$ws.Range("S11").Value = 0.003105590062111801
$ws.Range("G12").Value = 0.8585858585858586
$ws.Range("J12").Value = 0.1060606060606061
$ws.Range("L12").Value = 0.02525252525252525
$ws.Range("S12").Value = 0.0101010101010101
$ws.Range("G13").Value = 0.7924528301886793
$ws.Range("J13").Value = 0.2075471698113208
$ws.Range("F15").Value = 0.04072398190045249
$ws.Range("H15").Value = 0.1855203619909502
$ws.Range("I15").Value = 0.09954751131221719
$ws.Range("J15").Value = 0.3755656108597285
$ws.Range("K15").Value = 0.04977375565610859
$ws.Range("M15").Value = 0.01357466063348416
$ws.Range("O15").Value = 0.04072398190045249
$ws.Range("S15").Value = 0.1945701357466063
$ws.Range("F16").Value = 0.01111111111111111
$ws.Range("H16").Value = 0.1944444444444444
$ws.Range("I16").Value = 0.15
$ws.Range("J16").Value = 0.3833333333333334
$ws.Range("K16").Value = 0.08888888888888889
$ws.Range("M16").Value = 0.005555555555555556
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.02973977695167286
$ws.Range("H17").Value = 0.1672862453531599
$ws.Range("I17").Value = 0.1189591078066914
$ws.Range("J17").Value = 0.4144981412639405
$ws.Range("K17").Value = 0.09665427509293681
$ws.Range("M17").Value = 0.0241635687732342
$ws.Range("O17").Value = 0.06319702602230483
$ws.Range("S17").Value = 0.08550185873605948
$ws.Range("F18").Value = 0.02590673575129534
$ws.Range("H18").Value = 0.1917098445595855
$ws.Range("I18").Value = 0.1606217616580311
$ws.Range("J18").Value = 0.3367875647668394
$ws.Range("K18").Value = 0.1088082901554404
$ws.Range("M18").Value = 0.02590673575129534
$ws.Range("O18").Value = 0.07253886010362694
$ws.Range("S18").Value = 0.07772020725388601
$ws.Range("F19").Value = 0.0130246020260492
$ws.Range("H19").Value = 0.2040520984081042
$ws.Range("I19").Value = 0.130246020260492
$ws.Range("J19").Value = 0.357452966714906
$ws.Range("K19").Value = 0.1143270622286541
$ws.Range("M19").Value = 0.02243125904486252
$ws.Range("N19").Value = 0.001447178002894356
$ws.Range("O19").Value = 0.0723589001447178
$ws.Range("S19").Value = 0.08465991316931983
